# Add a new "Web2" page dependency row in the ManRelations sheet.
# This inserts a new row at position 18 (shifting the existing rows 18-21
# down to 19-22) and fills it with the new entry:
#   Enable=YES, start=www.magenta.ca|order, type=DependOf, end=MGTAWeb2,
#   weight=1, category=Availability, owner=IA

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(18).Insert() | Out-Null

$ws.Range("A18").Value = "YES"
$ws.Range("B18").Value = "www.magenta.ca|order"
$ws.Range("C18").Value = "DependOf"
$ws.Range("D18").Value = "MGTAWeb2"
$ws.Range("E18").Value = 1
$ws.Range("F18").Value = "Availability"
$ws.Range("G18").Value = "IA"

$ws.Range("B18").Select() | Out-Null
